# Open the active workbook / sheet (Plan1 - "Agenda de projeto")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Status" column (F) for rows 12-18 moves from "?" to "!"
# (i.e. from "em aberto" to "concluída", per the legend in the F3 comment)
$ws.Range("F12:F18").Value = "!"

# Move the current selection to A21 (and, since the view is no longer
# scrolled to keep A4 pinned at the top, drop the old frozen/topLeft
# viewport by just selecting the new cell)
$ws.Range("A21").Select()
